$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by the edit: A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values of rows 2, 4 and 5 before overwriting anything,
# since the edit performs a cyclic rotation of data among these rows:
#   new row2 = old row4
#   new row4 = old row5
#   new row5 = old row2
$row2 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row5[$col]
    $ws.Range("${col}5").Value = $row2[$col]
}
